# cats.xlsx update: rename lat/long column headers and refresh the saved
# view state (scroll back to the top of the sheet, reselect D2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers: "CenterLat"/"CenterLong" -> "Lat"/"Lon"
$ws.Range("D1").Value = "Lat"
$ws.Range("E1").Value = "Lon"

# Reset the saved scroll position (was topLeftCell="A10") back to the
# top-left of the sheet, and move the selection from H2 to D2.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D2").Select()
